$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 227 - 2021-04-15 (15 marzo per commit message numbering)
$ws.Cells.Item(226, 1).Copy($ws.Cells.Item(227, 1))
$ws.Cells.Item(227, 1).Value = 44301
$ws.Cells.Item(227, 2).Value = 2
$ws.Cells.Item(227, 3).Value = 25
$ws.Cells.Item(227, 4).Value = 294.5681630729351

# Row 228 - 2021-04-16
$ws.Cells.Item(226, 1).Copy($ws.Cells.Item(228, 1))
$ws.Cells.Item(228, 1).Value = 44302
$ws.Cells.Item(228, 2).Value = 8
$ws.Cells.Item(228, 3).Value = 22
$ws.Cells.Item(228, 4).Value = 259.2199835041828

# Row 229 - 2021-04-17
$ws.Cells.Item(226, 1).Copy($ws.Cells.Item(229, 1))
$ws.Cells.Item(229, 1).Value = 44303
$ws.Cells.Item(229, 2).Value = 7
$ws.Cells.Item(229, 3).Value = 24
$ws.Cells.Item(229, 4).Value = 282.7854365500177
